{"js": "// Update the table caption/title text and the ATT estimate figures\n// (caption changed from \"Table 3 ... with Controls\" to \"Table 2 ...\"\n// and the numeric results were regenerated without controls).\n\n// Helper: replace the first search hit of `oldText` scoped to a given\n// search-result collection with `newText`, preserving the run/paragraph\n// formatting (search-range insertText only swaps the text run).\nasync function replaceOnce(searchScope, oldText, newText) {\n  const results = searchScope.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1) Caption paragraph text.\nawait replaceOnce(\n  body,\n  \"Table 3. Aggregated and Cohort-Specific ATT Estimates with Controls\",\n  \"Table 2. Aggregated and Cohort-Specific ATT Estimates\"\n);\n\n// 2) Table caption property (w:tblCaption on w:tblPr).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\ntable.title = \"Table 2. Aggregated and Cohort-Specific ATT Estimates\";\nawait context.sync();\n\n// 3) Numeric results. \"-9.93\" appears in both row 1 and row 2 with\n// different replacement targets, so those two are scoped to their own\n// cell; the rest are unique across the whole document body.\nconst row1Att = table.getCell(1, 1); // \"Overall ATT (Group aggregation)\" / ATT_pp\nawait replaceOnce(row1Att.body, \"-9.93\", \"-1.97\");\n\nconst row2Att = table.getCell(2, 1); // \"Overall ATT (Dynamic aggregation)\" / ATT_pp\nawait replaceOnce(row2Att.body, \"-9.93\", \"-2.11\");\n\nconst otherReplacements = [\n  [\"1.79\", \"1.51\"],\n  [\"-13.44\", \"-4.94\"],\n  [\"-6.42\", \"0.99\"],\n  [\"2.23\", \"1.41\"],\n  [\"-14.31\", \"-4.88\"],\n  [\"-5.56\", \"0.66\"],\n  [\"-12.51\", \"-2.92\"],\n  [\"3.41\", \"1.63\"],\n  [\"-19.19\", \"-6.12\"],\n  [\"-5.83\", \"0.27\"],\n  [\"-7.79\", \"-1.18\"],\n  [\"1.82\", \"2.46\"],\n  [\"-11.34\", \"-6.00\"],\n  [\"-4.23\", \"3.63\"],\n];\n\nfor (const [oldText, newText] of otherReplacements) {\n  await replaceOnce(body, oldText, newText);\n}\n", "ps1": "# Update the table caption/title and the ATT estimate figures.\n# (caption changed from \"Table 3 ... with Controls\" to \"Table 2 ...\"\n# and the numeric results were regenerated without controls.)\n\n$d = $word.ActiveDocument\n\n# 1) Caption paragraph text (first paragraph in the document body).\n$d.Paragraphs.Item(1).Range.Text = \"Table 2. Aggregated and Cohort-Specific ATT Estimates\"\n\n# 2) Table caption property (w:tblCaption on w:tblPr).\n$tbl = $d.Tables.Item(1)\n$tbl.Title = \"Table 2. Aggregated and Cohort-Specific ATT Estimates\"\n\n# 3) Numeric results, by row/column (row 1 is the header row).\n#    Columns: 1=Row, 2=ATT_pp, 3=Std_Error, 4=CI_Lower_95, 5=CI_Upper_95\n\n# Row 2: \"Overall ATT (Group aggregation)\"\n$tbl.Cell(2, 2).Range.Text = \"-1.97\"\n$tbl.Cell(2, 3).Range.Text = \"1.51\"\n$tbl.Cell(2, 4).Range.Text = \"-4.94\"\n$tbl.Cell(2, 5).Range.Text = \"0.99\"\n\n# Row 3: \"Overall ATT (Dynamic aggregation)\"\n$tbl.Cell(3, 2).Range.Text = \"-2.11\"\n$tbl.Cell(3, 3).Range.Text = \"1.41\"\n$tbl.Cell(3, 4).Range.Text = \"-4.88\"\n$tbl.Cell(3, 5).Range.Text = \"0.66\"\n\n# Row 4: \"Cohort 2014\"\n$tbl.Cell(4, 2).Range.Text = \"-2.92\"\n$tbl.Cell(4, 3).Range.Text = \"1.63\"\n$tbl.Cell(4, 4).Range.Text = \"-6.12\"\n$tbl.Cell(4, 5).Range.Text = \"0.27\"\n\n# Row 5: \"Cohort 2015\"\n$tbl.Cell(5, 2).Range.Text = \"-1.18\"\n$tbl.Cell(5, 3).Range.Text = \"2.46\"\n$tbl.Cell(5, 4).Range.Text = \"-6.00\"\n$tbl.Cell(5, 5).Range.Text = \"3.63\"\n"}
